$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 411 (old rows 411:422 shift down to 414:425).
$ws.Rows("411:413").Insert()

# New weekly price records (Ají / Hortaliza, Mercado Mayorista Lo Valledor de Santiago).
$newRows = @(
    @{ Row = 411; A = 6; B = "Mercado Mayorista Lo Valledor de Santiago"; C = "Metropolitana"; D = 44448; E = 13; F = 100112021; G = "Ají"; H = "Americana (o)"; I = "Primera"; J = 130;  K = 90000; L = 95000; M = 93077; N = "`$/caja 25 kilos"; O = "Provincia de Limarí"; P = 3723; Q = 25; R = "Hortaliza" },
    @{ Row = 412; A = 6; B = "Mercado Mayorista Lo Valledor de Santiago"; C = "Metropolitana"; D = 44448; E = 13; F = 100112021; G = "Ají"; H = "Americana (o)"; I = "Segunda"; J = 40;   K = 85000; L = 85000; M = 85000; N = "`$/caja 25 kilos"; O = "Provincia de Limarí"; P = 3400; Q = 25; R = "Hortaliza" },
    @{ Row = 413; A = 6; B = "Mercado Mayorista Lo Valledor de Santiago"; C = "Metropolitana"; D = 44448; E = 13; F = 100112021; G = "Ají"; H = "Inferno";        I = "Primera"; J = 40;   K = 35000; L = 40000; M = 38000; N = "`$/caja 12 kilos"; O = "Provincia de Limarí"; P = 3167; Q = 12; R = "Hortaliza" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
}
